# Update the "Correspond Handoff Datetime" (D3) and
# "Correspond Handback DateTime" (G3) values on the zh-cn and de-de
# sheets to reflect the regenerated report timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 02:18:32"
$wsZhCn.Range("G3").Value = "2016-01-18 02:19:18"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 02:18:44"
$wsDeDe.Range("G3").Value = "2016-01-18 02:19:38"
